# Refresh the cryptos price/volume table with the latest scrape.
# Only the numeric "Price" (D) and "Volume(1h)" (E) columns move for most
# rows; three pairs of rows (33/34, 42/43, 49/50) also swapped rank order,
# so their Coin name (B) and Link (C) cells are rewritten too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> { column letter -> new text value }
$updates = @{
    2  = @{ D = "28.405.48"; E = "  -0.24%  " }
    3  = @{ D = "1.797.16";  E = "  -0.83%  " }
    4  = @{ D = "1.005";     E = "  +0.54%  " }
    5  = @{ D = "316.70";    E = "  -0.08%  " }
    6  = @{ D = "1.003";     E = "  +0.34%  " }
    7  = @{ D = "0.5405";    E = "  -2.38%  " }
    8  = @{ D = "0.3781";    E = "  -1.92%  " }
    9  = @{ D = "0.07483";   E = "  -1.55%  " }
    10 = @{ D = "42.01";     E = "  -2.18%  " }
    11 = @{ D = "1.109";     E = "  -2.04%  " }
    12 = @{ D = "1.005";     E = "  +0.59%  " }
    13 = @{ D = "20.61";     E = "  -3.24%  " }
    14 = @{ D = "6.150";     E = "  -1.11%  " }
    15 = @{ D = "7.306";     E = "  -0.56%  " }
    16 = @{ D = "1.794.07";  E = "  -0.65%  " }
    17 = @{ D = "89.60";     E = "  -1.74%  " }
    18 = @{ D = "0.00001065"; E = "  -0.87%  " }
    19 = @{ D = "0.06524";   E = "  +0.88%  " }
    20 = @{ D = "17.47";     E = "  +0.92%  " }
    21 = @{ E = "  +0.17%  " }
    22 = @{ E = "  -0.98%  " }
    23 = @{ D = "28.471.70"; E = "  -0.02%  " }
    24 = @{ D = "11.12";     E = "  -1.69%  " }
    25 = @{ D = "2.083";     E = "  -1.31%  " }
    26 = @{ D = "159.75";    E = "  +2.27%  " }
    27 = @{ D = "20.49" }
    28 = @{ D = "1.999.29";  E = "  -0.73%  " }
    29 = @{ D = "2.319";     E = "  -4.93%  " }
    30 = @{ D = "123.01";    E = "  -0.88%  " }
    31 = @{ D = "1.114";     E = "  -5.48%  " }
    32 = @{ D = "0.1055";    E = "  +1.88%  " }
    33 = @{ B = "Filecoin";          C = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil";          D = "5.615"; E = "  -2.69%  " }
    34 = @{ B = "HuobiToken";        C = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht";     D = "3.661"; E = "  +0.19%  " }
    35 = @{ D = "0.06504";   E = "  +3.52%  " }
    36 = @{ D = "0.2248";    E = "  -1.98%  " }
    37 = @{ D = "0.02290";   E = "  -1.89%  " }
    38 = @{ D = "8.596";     E = "  -3.70%  " }
    39 = @{ D = "5.020";     E = "  -0.64%  " }
    40 = @{ D = "11.22";     E = "  -3.61%  " }
    41 = @{ D = "0.6194";    E = "  -3.33%  " }
    42 = @{ B = "WEMIXTOKEN";        C = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix";      D = "1.454"; E = "  +5.09%  " }
    43 = @{ B = "TrustWalletToken";  C = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt";  D = "1.191"; E = "  +1.69%  " }
    44 = @{ D = "1.002";     E = "  +0.29%  " }
    45 = @{ D = "13.25";     E = "  -1.78%  " }
    46 = @{ D = "3.688";     E = "  +0.18%  " }
    47 = @{ D = "0.5822";    E = "  -3.04%  " }
    48 = @{ D = "126.62";    E = "  +1.78%  " }
    49 = @{ B = "EOS";               C = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos";           D = "1.211"; E = "  +5.27%  " }
    50 = @{ B = "NEARProtocol";      C = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near";     D = "1.952"; E = "  -1.04%  " }
    51 = @{ D = "0.06892";   E = "  -0.71%  " }
}

foreach ($row in $updates.Keys) {
    $rowData = $updates[$row]
    foreach ($col in $rowData.Keys) {
        $cell = $ws.Range("$col$row")
        $value = $rowData[$col]

        if ($col -eq "D" -and $value -match '^-?[0-9]*\.?[0-9]+$') {
            # These "Price" cells are plain text in the workbook (several
            # legitimately look like numbers, e.g. "1.005" or "89.60").
            # Force text entry so Excel doesn't silently coerce them to a
            # number and drop the formatting (trailing zeros, etc.), then
            # restore the default style so no stray number format lingers
            # on the cell.
            $cell.NumberFormat = "@"
            $cell.Value = $value
            $cell.Style = "Normal"
        } else {
            $cell.Value = $value
        }
    }
}
